# tvdcc_report.xlsx update
# "updated reports with complexity=4 to show show complex units"
#
# 1. Summary sheet: lower the "complexity" parameter from 10 to 4, and
#    refresh the processing start/end timestamps + elapsed seconds that a
#    fresh run with the new parameter produced.
# 2. PLSQLUnits sheet: with complexity=4 the run now reports 7 PL/SQL units
#    that exceed the threshold (previously none) - append those rows.
# 3. Files sheet: refresh the per-file "time in seconds" column (R) with the
#    numbers produced by the re-run.
# 4. Workbook-level defined name UnitsTable must grow to cover the newly
#    added PLSQLUnits rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

# complexity parameter: 10 -> 4
$summary.Range("B21").Value = "4"

# processing start / end time (serial date-times) and elapsed seconds
$summary.Range("B36").Value = 44806.49607638889
$summary.Range("B37").Value = 44806.496145833335
$summary.Range("B38").Value = 5.768

# ---------------------------------------------------------------------
# 2. PLSQLUnits sheet - append the newly reported complex units
# ---------------------------------------------------------------------
$units = $wb.Worksheets.Item("PLSQLUnits")

# columns: A File name, B PL/SQL Unit, C Line, D # Lines, E # Comment lines,
# F # Blank lines, G # Net lines, H # Stmts, I Cyclomatic complexity,
# J Halstead volume, K MI
$unitRows = @(
    @('guidelines/guideline_1040_04.sql', 'AnonymousPlsqlBlock',       10, 31, 1, 6, 31, 14, 5, 349.0333754971396,  93.3400232177323),
    @('guidelines/guideline_4370_45.sql', 'AnonymousPlsqlBlock',       14, 29, 0, 3, 29, 13, 5, 411.1982937621106,  82.40588543809453),
    @('guidelines/guideline_4310_39.sql', 'my_package.password_check', 26, 21, 0, 3, 21, 10, 5, 491.54240635418904, 86.13827722867563),
    @('guidelines/guideline_4310_39.sql', 'my_package.password_check', 63, 21, 0, 3, 21, 10, 5, 491.54240635418904, 86.13827722867563),
    @('guidelines/guideline_4320_40.sql', 'AnonymousPlsqlBlock',       12, 18, 0, 3, 18,  9, 5, 288.85263754543286, 91.06596991130587),
    @('guidelines/guideline_4320_40.sql', 'AnonymousPlsqlBlock',       41, 25, 0, 3, 25,  9, 5, 346.1295543881475,  85.46485860912527),
    @('guidelines/guideline_4370_45.sql', 'AnonymousPlsqlBlock',       54, 25, 0, 3, 25,  9, 5, 346.1295543881475,  85.46485860912527)
)

$startRow = 2
for ($i = 0; $i -lt $unitRows.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $unitRows[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $units.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
    # Cyclomatic complexity / Halstead volume / MI columns use an integer
    # display format, matching the rest of the report's metric columns.
    $units.Range("I" + $r + ":K" + $r).NumberFormat = "0"
}

# ---------------------------------------------------------------------
# 3. Files sheet - refreshed "time in seconds" (column R) values
# ---------------------------------------------------------------------
$files = $wb.Worksheets.Item("Files")

$fileTimes = @{
    2 = 0.008; 3 = 0.008; 4 = 0.034; 5 = 0.017; 6 = 0.067; 7 = 0.06; 8 = 0.011;
    9 = 0.013; 10 = 0.011; 11 = 0.006; 12 = 0.013; 14 = 0.009; 15 = 0.011;
    17 = 0.011; 18 = 0.009; 19 = 0.008; 20 = 0.057; 21 = 0.009; 24 = 0.007;
    25 = 0.013; 26 = 0.013; 27 = 0.048; 29 = 0.007; 31 = 0.007; 33 = 0.01;
    34 = 0.012; 35 = 0.02; 36 = 0.008; 37 = 0.008; 38 = 0.01; 39 = 0.01;
    41 = 0.028; 42 = 0.033; 43 = 0.052; 44 = 0.012; 45 = 0.057; 47 = 0.059;
    48 = 0.012; 49 = 0.01; 50 = 0.008; 51 = 0.054; 52 = 0.048; 53 = 0.01;
    54 = 0.008; 55 = 0.006; 56 = 0.017; 57 = 0.009; 59 = 0.005; 60 = 0.007;
    61 = 0.012; 62 = 0.44; 63 = 0.014; 64 = 0.028; 65 = 0.014; 67 = 0.021;
    68 = 0.055; 69 = 0.009; 70 = 0.013; 71 = 0.015; 72 = 0.014; 73 = 0.011;
    74 = 0.007; 78 = 0.055; 79 = 0.009; 80 = 0.055; 81 = 0.053; 82 = 0.015;
    83 = 0.007; 84 = 0.011; 85 = 0.013; 86 = 0.017; 87 = 0.062; 88 = 0.009;
    90 = 0.006; 91 = 0.008; 94 = 0.007; 96 = 0.01; 97 = 0.009; 98 = 0.058;
    99 = 0.008; 102 = 0.006; 105 = 0.004; 106 = 0.006; 107 = 0.008;
    108 = 1.935; 110 = 0.008; 116 = 0.01; 118 = 0.008; 119 = 0.016;
    120 = 0.057; 121 = 0.054; 122 = 0.006
}

foreach ($row in $fileTimes.Keys) {
    $files.Range("R" + $row).Value = $fileTimes[$row]
}

# ---------------------------------------------------------------------
# 4. Workbook-level defined name - grow UnitsTable to cover the new rows
# ---------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "UnitsTable") {
        $n.RefersTo = "=PLSQLUnits!`$A`$1:`$H`$8"
    }
}
